# Refresh the cryptocurrency Price (column D) and Volume(1h) (column E)
# figures on the active worksheet to the latest scraped snapshot.
#
# The sheet stores these as plain text (not numbers) so that values such
# as "0.600" or "63.799.66" keep their exact original formatting. Excel's
# Range.Value setter auto-detects numeric-looking strings and silently
# converts them to real numbers, which would both reformat the text
# (dropping trailing zeros, etc.) and change the cell's stored type. To
# avoid that, any replacement Price that looks like a plain number is
# written while the cell is temporarily forced to Text format, and the
# formatting is cleared again immediately afterwards so the cell ends up
# with no explicit style applied (matching the rest of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "63.828.07"; DNumeric = $false; E = "  -0.74%  " },
    @{ Row = 3; D = "2.746.92"; DNumeric = $false; E = "  -0.86%  " },
    @{ Row = 4; D = $null; DNumeric = $false; E = "  +0.19%  " },
    @{ Row = 5; D = "572.63"; DNumeric = $true; E = "  -1.44%  " },
    @{ Row = 6; D = "157.17"; DNumeric = $true; E = "  +1.06%  " },
    @{ Row = 7; D = $null; DNumeric = $false; E = "  +0.22%  " },
    @{ Row = 8; D = "0.600"; DNumeric = $true; E = "  -1.42%  " },
    @{ Row = 9; D = $null; DNumeric = $false; E = "  -3.39%  " },
    @{ Row = 10; D = $null; DNumeric = $false; E = "  -0.22%  " },
    @{ Row = 12; D = "5.56"; DNumeric = $true; E = "  -17.53%  " },
    @{ Row = 13; D = "3.232.41"; DNumeric = $false; E = "  -0.51%  " },
    @{ Row = 14; D = "26.36"; DNumeric = $true; E = "  -1.75%  " },
    @{ Row = 15; D = "63.514.92"; DNumeric = $false; E = "  -0.59%  " },
    @{ Row = 16; D = $null; DNumeric = $false; E = "  -2.50%  " },
    @{ Row = 17; D = "2.750.15"; DNumeric = $false; E = "  -0.64%  " },
    @{ Row = 18; D = "12.11"; DNumeric = $true; E = "  +0.87%  " },
    @{ Row = 19; D = $null; DNumeric = $false; E = "  -2.18%  " },
    @{ Row = 20; D = "353.87"; DNumeric = $true; E = "  -2.37%  " },
    @{ Row = 21; D = "6.71"; DNumeric = $true; E = "  -4.23%  " },
    @{ Row = 22; D = $null; DNumeric = $false; E = "  +0.63%  " },
    @{ Row = 23; D = $null; DNumeric = $false; E = "  -0.15%  " },
    @{ Row = 24; D = "65.02"; DNumeric = $true; E = "  -2.14%  " },
    @{ Row = 25; D = $null; DNumeric = $false; E = "  -1.36%  " },
    @{ Row = 26; D = $null; DNumeric = $false; E = "  -0.01%  " },
    @{ Row = 27; D = "8.36"; DNumeric = $true; E = "  -2.37%  " },
    @{ Row = 28; D = "0.0₃0903"; DNumeric = $false; E = "  -0.29%  " },
    @{ Row = 29; D = $null; DNumeric = $false; E = "  -3.90%  " },
    @{ Row = 30; D = "6.93"; DNumeric = $true; E = "  -2.91%  " },
    @{ Row = 31; D = "168.97"; DNumeric = $true; E = "  -2.25%  " },
    @{ Row = 32; D = $null; DNumeric = $false; E = "  -8.17%  " },
    @{ Row = 33; D = $null; DNumeric = $false; E = "  -2.13%  " },
    @{ Row = 34; D = $null; DNumeric = $false; E = "  +0.18%  " },
    @{ Row = 35; D = "4.84"; DNumeric = $true; E = "  -0.71%  " },
    @{ Row = 36; D = $null; DNumeric = $false; E = "  -1.13%  " },
    @{ Row = 37; D = "1.79"; DNumeric = $true; E = "  -2.74%  " },
    @{ Row = 38; D = "0.976"; DNumeric = $true; E = "  -3.81%  " },
    @{ Row = 39; D = "6.15"; DNumeric = $true; E = "  +4.58%  " },
    @{ Row = 40; D = $null; DNumeric = $false; E = "  -3.09%  " },
    @{ Row = 41; D = "324.16"; DNumeric = $true; E = "  -6.31%  " },
    @{ Row = 42; D = "38.81"; DNumeric = $true; E = "  -1.42%  " },
    @{ Row = 43; D = "21.20"; DNumeric = $true; E = "  -3.87%  " },
    @{ Row = 44; D = "0.0583"; DNumeric = $true; E = "  -1.91%  " },
    @{ Row = 45; D = "21.22"; DNumeric = $true; E = "  -3.82%  " },
    @{ Row = 46; D = $null; DNumeric = $false; E = "  -1.85%  " },
    @{ Row = 47; D = "134.65"; DNumeric = $true; E = "  -2.22%  " },
    @{ Row = 48; D = $null; DNumeric = $false; E = "  -4.36%  " },
    @{ Row = 49; D = $null; DNumeric = $false; E = "  -1.08%  " },
    @{ Row = 50; D = $null; DNumeric = $false; E = "  -0.08%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        if ($u.DNumeric) {
            $dCell.NumberFormat = "@"
            $dCell.Value = $u.D
            $dCell.ClearFormats()
        } else {
            $dCell.Value = $u.D
        }
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
